# Economic Dashboard V1 - weekly data refresh (2025-11-27)
# Applies the source-diff's cell-level updates: new FRED pull dates, refreshed
# rolling-window deltas for Durable Goods / Non-Def-ex-Air orders, updated
# Treasury / credit-spread "as of" columns, and the highlight-fill toggle that
# marks which NFP/UR/LFPR/EPOP/AHE/AWH series were refreshed this pull vs.
# which bond-yield row is now "stale" (fill removed) / "fresh" (fill added).
#
# NOTE: this COM host does not reliably replay custom function calls that
# wrap Range/Copy/PasteSpecial, so every step below is inlined directly
# (loops over a literal address list are fine, user-defined functions are
# not used).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Stable reference cells (their own style never changes in this edit):
#   C3  -> style with NO fill  (plain "as of" date)
#   N47 -> style WITH yellow fill (highlighted "as of" date)

# ---------------------------------------------------------------------------
# 1) Remove highlight fill from the "as of" date cells that are no longer the
#    freshest pull (style 48 -> 47); the date value itself is unchanged.
#    Reuse the existing no-fill style by copying format from C3.
# ---------------------------------------------------------------------------
foreach ($addr in @("N3","N4","N6","N7","N8","N9","N15","N33","N34","N35","N36")) {
    $ws.Range("C3").Copy()
    $ws.Range($addr).PasteSpecial(-4122)
}

# ---------------------------------------------------------------------------
# 2) Add highlight fill to "as of" date cells that are newly the freshest pull
#    (style 47 -> 48). Reuse the existing highlighted style by copying format
#    from N47.
# ---------------------------------------------------------------------------
foreach ($addr in @("C7","C28","C29","C30","C31")) {
    $ws.Range("N47").Copy()
    $ws.Range($addr).PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3) GDPNow (row 7): refreshed nowcast value.
# ---------------------------------------------------------------------------
$ws.Range("F7").Value2 = 0.3862508614748448

# ---------------------------------------------------------------------------
# 4) Durable Goods Orders, M/M (row 28) - rolling window shifts + new "as of".
# ---------------------------------------------------------------------------
$ws.Range("C28").Value2 = 45901
$ws.Range("F28").Value2 = 0.004840745558111426
$ws.Range("G28").Value2 = 0.03006283164814283
$ws.Range("H28").Value2 = -0.02799901206372835
$ws.Range("I28").Value2 = -0.09389977010425232
$ws.Range("J28").Value2 = 0.165430902790715

# ---------------------------------------------------------------------------
# 5) Durable Goods Orders, Y/Y (row 29) - rolling window shifts + new "as of".
# ---------------------------------------------------------------------------
$ws.Range("C29").Value2 = 45901
$ws.Range("F29").Value2 = 0.07243160813624692
$ws.Range("G29").Value2 = 0.07662644950780544
$ws.Range("H29").Value2 = 0.03341358778313566
$ws.Range("I29").Value2 = 0.1089645997552716
$ws.Range("J29").Value2 = 0.1995614111543812

# 5yr,5yr Forward (row 29, right block): new "as of" date + shifted daily series.
$ws.Range("N29").Value2 = 45987
$ws.Range("R29").Value2 = 2.17
$ws.Range("S29").Value2 = 2.16
$ws.Range("U29").Value2 = $null

# ---------------------------------------------------------------------------
# 6) Dur Orders Non Def x Aircraft, M/M (row 30) - rolling window + "as of".
# ---------------------------------------------------------------------------
$ws.Range("C30").Value2 = 45901
$ws.Range("F30").Value2 = 0.0006782625461356773
$ws.Range("G30").Value2 = 0.01909075917433611
$ws.Range("H30").Value2 = -0.02404555711932721
$ws.Range("I30").Value2 = -0.09442194506291901
$ws.Range("J30").Value2 = 0.1570751450479186

# 10yr TIPS (row 30, right block): new "as of" date + shifted daily series.
$ws.Range("N30").Value2 = 45987
$ws.Range("Q30").Value2 = 2.23
$ws.Range("R30").Value2 = 2.22
$ws.Range("S30").Value2 = 2.23
$ws.Range("U30").Value2 = $null

# ---------------------------------------------------------------------------
# 7) Dur Orders Non Def x Aircraft, Y/Y (row 31) - rolling window + "as of".
# ---------------------------------------------------------------------------
$ws.Range("C31").Value2 = 45901
$ws.Range("F31").Value2 = 0.06450892203111723
$ws.Range("G31").Value2 = 0.06672542970471573
$ws.Range("H31").Value2 = 0.0329297153895499
$ws.Range("I31").Value2 = 0.1029410098461701
$ws.Range("J31").Value2 = 0.1925388202490265

# ---------------------------------------------------------------------------
# 8) FFR / 2y / 5y / 10y UST, 30y Mortgage, BAA "as of" dates roll forward;
#    their daily-value columns (Q..U) shift right as a new reading lands.
# ---------------------------------------------------------------------------
$ws.Range("N47").Value2 = 45986

$ws.Range("N48").Value2 = 45986
$ws.Range("Q48").Value2 = 3.43
$ws.Range("R48").Value2 = 3.46
$ws.Range("T48").Value2 = $null
$ws.Range("U48").Value2 = 3.51

$ws.Range("N49").Value2 = 45986
$ws.Range("Q49").Value2 = 3.55
$ws.Range("R49").Value2 = 3.61
$ws.Range("T49").Value2 = $null
$ws.Range("U49").Value2 = 3.62

$ws.Range("N50").Value2 = 45986
$ws.Range("Q50").Value2 = 4.01
$ws.Range("R50").Value2 = 4.04
$ws.Range("T50").Value2 = $null
$ws.Range("U50").Value2 = 4.06

# 30y Mortgage (weekly series) - "as of" date + full shift of the 5 columns.
$ws.Range("N51").Value2 = 45985
$ws.Range("Q51").Value2 = 6.23
$ws.Range("R51").Value2 = 6.26
$ws.Range("S51").Value2 = 6.24
$ws.Range("T51").Value2 = 6.22
$ws.Range("U51").Value2 = 6.17

# BAA corporate bond yield.
$ws.Range("N52").Value2 = 45986
$ws.Range("Q52").Value2 = 5.8
$ws.Range("R52").Value2 = 5.84
$ws.Range("T52").Value2 = $null
$ws.Range("U52").Value2 = 5.88
